$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4 and 5 (companies no longer present in the dataset)
$ws.Rows("4:5").Delete()

# Row 2 updates
$ws.Range("B2").Value = "'1"
$ws.Range("D2").Value = -0.475
$ws.Range("E2").Value = -0.133
$ws.Range("G2").Value = 1.003703703703704
$ws.Range("H2").Value = 1.003703703703704
$ws.Range("I2").Value = 0.5900150401356185
$ws.Range("J2").Value = 0.5810962081335684
$ws.Range("K2").Value = 2.54
$ws.Range("L2").Value = 0.4703703703703703
$ws.Range("M2").Value = 3.19437
$ws.Range("N2").Value = 0.05202557003257329
$ws.Range("O2").Value = 1.257625984251969
$ws.Range("P2").Value = 1.64437
$ws.Range("Q2").Value = 0.02678127035830619
$ws.Range("R2").Value = 0.6473897637795276
$ws.Range("S2").Value = 1.55
$ws.Range("T2").Value = 0.4852286992427302
$ws.Range("U2").Value = 1.46
$ws.Range("V2").Value = 0.0237785016286645
$ws.Range("W2").Value = 0.03227445997458704
$ws.Range("X2").Value = 0.03625635443668625
$ws.Range("Y2").Value = -0.003981894462099213
$ws.Range("Z2").Value = 0.09420862275963876
$ws.Range("AA2").Value = 0.05474427345911188
$ws.Range("AB2").Value = 0.03517703359074804
$ws.Range("AC2").Value = 0.01956723986836384
$ws.Range("AD2").Value = 3.01
$ws.Range("AE2").Value = 0.1795939163383019
$ws.Range("AF2").Value = 3.189593916338302
$ws.Range("AG2").Value = 1.729593916338302
$ws.Range("AH2").Value = 0.04938247359891631
$ws.Range("AI2").Value = 0.03967670143548332
$ws.Range("AJ2").Value = 0.02739751373389831
$ws.Range("AK2").Value = 0.02191312320916779
$ws.Range("AL2").Value = 0.483
$ws.Range("AM2").Value = 0.464
$ws.Range("AN2").Value = 0.9255842558425584
$ws.Range("AO2").Value = 6.418219461697723
$ws.Range("AP2").Value = 0.531855447828506
$ws.Range("AQ2").Value = 6.681034482758621

# Row 3 updates
$ws.Range("B3").Value = 'Value8 N.V. (ENXTAM:VALUE)'
$ws.Range("D3").Value = -0.475
$ws.Range("E3").Value = -0.133
$ws.Range("G3").Value = 1.003703703703704
$ws.Range("H3").Value = 1.003703703703704
$ws.Range("I3").Value = 0.5900150401356185
$ws.Range("J3").Value = 0.5810962081335684
$ws.Range("K3").Value = 2.54
$ws.Range("L3").Value = 0.4703703703703703
$ws.Range("M3").Value = 3.19437
$ws.Range("N3").Value = 0.05202557003257329
$ws.Range("O3").Value = 1.257625984251969
$ws.Range("P3").Value = 1.64437
$ws.Range("Q3").Value = 0.02678127035830619
$ws.Range("R3").Value = 0.6473897637795276
$ws.Range("S3").Value = 1.55
$ws.Range("T3").Value = 0.4852286992427302
$ws.Range("U3").Value = 1.46
$ws.Range("V3").Value = 0.0237785016286645
$ws.Range("W3").Value = 0.03227445997458704
$ws.Range("X3").Value = 0.03625635443668625
$ws.Range("Y3").Value = -0.003981894462099213
$ws.Range("Z3").Value = 0.09420862275963876
$ws.Range("AA3").Value = 0.05474427345911188
$ws.Range("AB3").Value = 0.03517703359074804
$ws.Range("AC3").Value = 0.01956723986836384
$ws.Range("AD3").Value = 3.01
$ws.Range("AE3").Value = 0.1795939163383019
$ws.Range("AF3").Value = 3.189593916338302
$ws.Range("AG3").Value = 1.729593916338302
$ws.Range("AH3").Value = 0.04938247359891631
$ws.Range("AI3").Value = 0.03967670143548332
$ws.Range("AJ3").Value = 0.02739751373389831
$ws.Range("AK3").Value = 0.02191312320916779
$ws.Range("AL3").Value = 0.483
$ws.Range("AM3").Value = 0.464
$ws.Range("AN3").Value = 0.9255842558425584
$ws.Range("AO3").Value = 6.418219461697723
$ws.Range("AP3").Value = 0.531855447828506
$ws.Range("AQ3").Value = 6.681034482758621
